$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "Save" in H1, formatted like the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Row -> Save value (1 for the top-4 highest "sum" games, 0 otherwise)
$saveValues = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 1; 8 = 0; 9 = 0; 10 = 0;
    11 = 0; 12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 0; 18 = 0; 19 = 0;
    20 = 0; 21 = 1; 22 = 0; 23 = 0; 24 = 0; 25 = 1; 26 = 0; 27 = 0; 28 = 1;
    29 = 0; 30 = 0; 31 = 0; 32 = 0; 33 = 0; 34 = 0; 35 = 0; 36 = 0; 37 = 0;
    38 = 0; 39 = 0; 40 = 0; 41 = 0; 42 = 0; 43 = 0; 44 = 0; 45 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
